# "add data and fix graphic"
# Adds new Client rows (38-46), fills in the employee rows (83-88) with
# generated credentials, tweaks the Admin(HR) row for Navid Panahi (89),
# renames the password column header, widens column F and updates the
# active selection / window view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Client" rows for ids 9220015041-9220015049 (rows 38-46) ---------
$clientRows = @(
    @{ Row = 38; First = "Mohsen";  Last = "Jamebozorg" },
    @{ Row = 39; First = "Sina";    Last = "Ahmadi" },
    @{ Row = 40; First = "Mobina";  Last = "Rajabi" },
    @{ Row = 41; First = "Hassan";  Last = "Solymani" },
    @{ Row = 42; First = "Matin";   Last = "Movahidi" },
    @{ Row = 43; First = "Kazem";   Last = "Housiny" },
    @{ Row = 44; First = "Rahele";  Last = "Asadi" },
    @{ Row = 45; First = "Amir";    Last = "Taheri" },
    @{ Row = 46; First = "Kokab";   Last = "Ahmadi" }
)

foreach ($r in $clientRows) {
    $ws.Cells.Item($r.Row, 2).Value = "Client"
    $ws.Cells.Item($r.Row, 3).Value = $r.First
    $ws.Cells.Item($r.Row, 4).Value = $r.Last
}

# --- Fill in "employee" rows 83-88 with names + generated passwords -------
$ws.Cells.Item(83, 2).Value = "employee"
$ws.Cells.Item(83, 3).Value = "Ahmad"

# --- Header: rename "password" -> "password(meli code)" -------------------
$ws.Cells.Item(1, 6).Value = "password(meli code)"

$ws.Cells.Item(83, 4).Value = "Shakiri"
$ws.Cells.Item(83, 5).Value = 12345
$ws.Cells.Item(83, 6).Value = 12345

$employeeRows = @(
    @{ Row = 84; First = "Muhammad"; Last = "Qorbanzade"; Code = 123456 },
    @{ Row = 85; First = "Shayan";   Last = "Abdolahy";   Code = 1234567 },
    @{ Row = 86; First = "Kiana";    Last = "Iravani";    Code = 22334455 },
    @{ Row = 87; First = "Mahdi";    Last = "Kushanmehr"; Code = 11225522 },
    @{ Row = 88; First = "Yalda";    Last = "Tahbaz";     Code = 778842 }
)

foreach ($r in $employeeRows) {
    $ws.Cells.Item($r.Row, 2).Value = "employee"
    $ws.Cells.Item($r.Row, 3).Value = $r.First
    $ws.Cells.Item($r.Row, 4).Value = $r.Last
    $ws.Cells.Item($r.Row, 5).Value = $r.Code
    $ws.Cells.Item($r.Row, 6).Value = $r.Code
}

# --- Fix Admin(HR) row 89 (Navid Panahi) phone/password values ------------
$ws.Cells.Item(89, 5).Value = 987654331
$ws.Cells.Item(89, 6).Value = 98766

# --- Column F is wider now to fit the longer header -----------------------
$ws.Columns("F").ColumnWidth = 17.5

# --- Update the active cell / selection ------------------------------------
$ws.Range("F14").Select()
